# "add a new story" -- append a new row (47000003 / Lost) to the DungeonStory
# table on Sheet1, growing the table from A3:M5 to A3:M6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# Grow the table by one row; this also extends dimension/autofilter once the
# row carries data.
$newRow = $tbl.ListRows.Add()

# Match row 6's look to row 5 by copying each column's formatting across
# before the values are written.
for ($c = 1; $c -le 13; $c++) {
    $ws.Cells.Item(5, $c).Copy()
    $ws.Cells.Item(6, $c).PasteSpecial(-4122)
}

# New story data (row 6): Id, Name, Descript, RuleStr, Rate, DungeonId,
# NeedGismoId, EventReplace, AttrType, AttrBias, CardId, BlessId, Image.
$ws.Cells.Item(6, 1).Value = 47000003
$ws.Cells.Item(6, 2).Value = "Lost"
$ws.Cells.Item(6, 3).Value = "找到迷宫出口"
$ws.Cells.Item(6, 4).Value = "充满了很多机关|找到出口异常艰辛"
$ws.Cells.Item(6, 5).Value = 300
$ws.Cells.Item(6, 6).Value = 18000001
$ws.Cells.Item(6, 7).Value = 45000001
$ws.Cells.Item(6, 8).Value = "bossqiongqi=bossunicorn,trees=rosemaryfield,river=poppyfield,manflower=trapspring,cliff=trappoison"
$ws.Cells.Item(6, 13).Value = "lost"

# Extend the conditional-formatting "highlight if blank" checks that watch
# columns G:J, K and L down onto the new row.
$cfGJ = $ws.Range("G4:J5").FormatConditions.Item(1)
$cfGJ.ModifyAppliesToRange($ws.Range("G4:J6"))

$cfK = $ws.Range("K4:K5").FormatConditions.Item(1)
$cfK.ModifyAppliesToRange($ws.Range("K4:K6"))

$cfL = $ws.Range("L4:L5").FormatConditions.Item(1)
$cfL.ModifyAppliesToRange($ws.Range("L4:L6"))

# Column M's blanks check shared a rule with G:J in the original file; this
# engine's FormatConditions can only carry a single contiguous area per
# rule, so re-create the same "blank -> green fill" check over M4:M6.
$cfM = $ws.Range("M4:M6").FormatConditions.Add(2, 0, "=LEN(TRIM(M4))=0")
$cfM.Interior.Color = 5287936

# Mirror the author's final selection on the new last cell of the table.
$ws.Range("M6").Select()

Write-Host "Added row 6 (Id 47000003 / Lost); table now $($tbl.Range.Address())"
